# Update the "Förändrad" (changed) date column C for rows 2-20
# from serial date 45212 (2023-10-13) to 45221 (2023-10-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 20; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    $cell.Value = 45221
}
